# Tester A filled in her results for the three checklist items in her
# column of the test matrix (E7:E9 - "Pass"/"Fail" values validated
# against the $H$3:$H$4 list).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

$ws.Range("E7").Value = "Fail"
$ws.Range("E8").Value = "Pass"
$ws.Range("E9").Value = "Pass"

# Leave the selection where the author's cursor ended up.
$ws.Range("E11").Select()
